# sua chiet khau cua sale phu va update chien luoc chay tinh luong theo gio

$wb = $excel.ActiveWorkbook

# --- Sheet "Đơn phụ phẫu 1": insert two new detail rows before the Tổng row ---
$ws2 = $wb.Worksheets.Item("Đơn phụ phẫu 1")

# Push the existing "Tổng" row (row 11) down to row 13 to make room for the
# two new service rows.
$ws2.Rows("11").Insert()
$ws2.Rows("11").Insert()

# Row 11 - new data row
$ws2.Cells.Item(11, 1).Value = "HD-LUXURY"
$ws2.Cells.Item(11, 2).Value = 575
# The "Ngày thực hiện" column holds date-like text (e.g. "07-21-2024") that
# must stay plain text, not get auto-parsed into a date serial number:
# format the cell as Text first, write the value, then drop back to the
# (unformatted) default style so the stored text isn't re-interpreted.
$ws2.Cells.Item(11, 3).NumberFormat = "@"
$ws2.Cells.Item(11, 3).Value = "07-21-2024"
$ws2.Cells.Item(11, 3).ClearFormats()
$ws2.Cells.Item(11, 4).Value = "CẦN THƠ"
$ws2.Cells.Item(11, 5).Value = "Nguyễn Thị Lan Anh"
$ws2.Cells.Item(11, 6).Value = "Cá nhân"
$ws2.Cells.Item(11, 7).Value = "Nâng mũi"
$ws2.Cells.Item(11, 8).Value = "Lâm Hoàng Phú"
$ws2.Cells.Item(11, 9).Value = 100000

# Row 12 - new data row
$ws2.Cells.Item(12, 1).Value = "HD-LUXURY"
$ws2.Cells.Item(12, 2).Value = 576
$ws2.Cells.Item(12, 3).NumberFormat = "@"
$ws2.Cells.Item(12, 3).Value = "07-21-2024"
$ws2.Cells.Item(12, 3).ClearFormats()
$ws2.Cells.Item(12, 4).Value = "CẦN THƠ"
$ws2.Cells.Item(12, 5).Value = "Nguyễn Thị Phương"
$ws2.Cells.Item(12, 6).Value = "Cá nhân"
$ws2.Cells.Item(12, 7).Value = "Cắt mí"
$ws2.Cells.Item(12, 8).Value = "Lâm Hoàng Phú"
$ws2.Cells.Item(12, 9).Value = 50000

# Row 13 - Tổng (totals) row. It was row 11 before the insert above and the
# row-shift carried its (already blank) C:H cells down with it untouched, so
# only the count and grand-total cells need new values here.
$ws2.Cells.Item(13, 2).Value = 11
$ws2.Cells.Item(13, 9).Value = 850000

# --- Sheet "Lương": update derived totals ---
$ws3 = $wb.Worksheets.Item("Lương")

$ws3.Cells.Item(9, 2).Value = 850000
$ws3.Cells.Item(28, 2).Value = 1065000
$ws3.Cells.Item(31, 2).Value = 3169514.285714285
